# Add a new "2023" data column (Q) to the 4.c.1 sheet, mirroring the
# existing column P (2022) formatting, one column to its right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column P's formatting (rows 3-7) into column Q so the new cells
# pick up the same styles (borders/number formats/fonts) as the rest of
# the table instead of the sheet's bare default column style.
$ws.Range("P3:P7").Copy()
$ws.Range("Q3:Q7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New 2023 figures, one column to the right of the 2022 ("P") column.
$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 93.7
$ws.Range("Q6").Value = 95.5
$ws.Range("Q7").Value = 97.1
